$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D12").Value = "“머신 러닝 교과서 3판”이 출간되었습니다."
$ws.Range("E12").Value = "https://tensorflow.blog/2021/03/23/%eb%a8%b8%ec%8b%a0-%eb%9f%ac%eb%8b%9d-%ea%b5%90%ea%b3%bc%ec%84%9c-3%ed%8c%90%ec%9d%b4-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

$ws.Range("D36").Value = "Transformer in Computer Vision"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/316"

$ws.Range("D37").Value = "[Paper Review] Time Series Anomaly Detection with Multiresolution Ensemble Decoding"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1451&mod=document&pageid=1"

$ws.Range("D39").Value = "3 Best (Often Better) Alternatives To Histograms, Avoid Binning Bias"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/3-Best-Often-Better-Alternatives-To-Histograms-Avoid-Binning-Bias-1"

$ws.Range("D51").Value = "[독후감] 생활코딩의 이고잉님의 강의를 다듬어서 출판한 <Do it! 지옥에서 온 문서 관리자 깃&깃허브 입문>"
$ws.Range("E51").Value = "https://bskyvision.com/1150"
